$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
#    We build it by copying the paragraph-mark/formatting of an existing
#    plain ("Normal" style) paragraph so the new paragraph does not carry
#    over the Heading1 paragraph style (and does not pick up any stray
#    rsid bookkeeping attributes that a direct Style re-assignment would
#    add).
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$normalSourcePara = $d.Paragraphs.Item(4)
$normalSourceRange = $d.Range($normalSourcePara.Range.Start, $normalSourcePara.Range.End)

$insertPos = $titlePara.Range.End
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.FormattedText = $normalSourceRange.FormattedText

$metaPara = $d.Paragraphs.Item(2)
$metaTextRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End - 1)
$metaTextRange.Text = "Meta description: Take a look at our review of Dice Dice Baby slot and play for free. Learn about gameplay, symbols, bonuses, and device compatibility."

# Bold just the "Meta description" label (first 16 characters, no colon).
$labelStart = $metaPara.Range.Start
$labelRange = $d.Range($labelStart, $labelStart + 16)
$labelRange.Bold = 1

# ---------------------------------------------------------------------------
# 2) Drop the duplicate bold "Play Dice Dice Baby Free - Review of Dice
#    Dice Baby Slot" paragraph that was left near the end of the document.
#    (Index shifted by +1 because of the paragraph inserted in step 1.)
# ---------------------------------------------------------------------------
$dupHeadingPara = $d.Paragraphs.Item(50)
$dupHeadingPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the italic "meta description" blurb at the very end with the
#    image-generation prompt text (still italic).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item(50)
$lastTextRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$lastTextRange.Text = 'Create a feature image for "Dice Dice Baby" in cartoon style featuring a happy Maya warrior with glasses. This image should have an adventurous and playful look, with a bright color scheme to highlight the excitement of the game. The warrior should be shown holding two dice in their hand, with a smile on their face, and glasses adding a touch of modernity to their traditional outfit. The background should be a jungle scene, with animated foliage and colorful flowers. Overall, the image should convey the fun and thrill of playing this game online.'

Write-Output "done"
